$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume table to latest scraped values.
# Numeric-looking price strings are prefixed with a literal apostrophe so Excel
# stores them as text (matching the original inlineStr cell type) instead of
# silently converting them to numbers and losing formatting/precision.

$ws.Range("D2").Value = "27.535.64"
$ws.Range("E2").Value = "  -0.73%  "
$ws.Range("D3").Value = "1.619.76"
$ws.Range("E3").Value = "  -1.73%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "'211.14"
$ws.Range("E5").Value = "  -0.89%  "
$ws.Range("E6").Value = "  -1.89%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "'22.74"
$ws.Range("E8").Value = "  -1.65%  "
$ws.Range("D9").Value = "'0.263"
$ws.Range("E9").Value = "  +2.26%  "
$ws.Range("E10").Value = "  +0.15%  "
$ws.Range("E11").Value = "  -0.46%  "
$ws.Range("D12").Value = "1.848.89"
$ws.Range("D13").Value = "1.627.22"
$ws.Range("E13").Value = "  -1.12%  "
$ws.Range("E14").Value = "  -0.50%  "
$ws.Range("E15").Value = "  -2.12%  "
$ws.Range("D16").Value = "'65.01"
$ws.Range("E16").Value = "  +1.31%  "
$ws.Range("D17").Value = "27.541.20"
$ws.Range("E17").Value = "  -0.54%  "
$ws.Range("D18").Value = "'230.06"
$ws.Range("E18").Value = "  -0.16%  "
$ws.Range("D19").Value = "0.0₃0720"
$ws.Range("E19").Value = "  -0.57%  "
$ws.Range("D20").Value = "'7.53"
$ws.Range("E20").Value = "  -1.81%  "
$ws.Range("E21").Value = "  +0.25%  "
$ws.Range("D22").Value = "'4.28"
$ws.Range("E22").Value = "  -0.72%  "
$ws.Range("E23").Value = "  +0.45%  "
$ws.Range("E24").Value = "  +6.42%  "
$ws.Range("D25").Value = "'149.55"
$ws.Range("E25").Value = "  +0.21%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").Value = "'6.82"
$ws.Range("E26").Value = "  -2.23%  "
$ws.Range("B27").Value = "Stellar"
$ws.Range("C27").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D27").Value = "'0.111"
$ws.Range("E27").Value = "  -1.46%  "
$ws.Range("E28").Value = "  +0.21%  "
$ws.Range("E29").Value = "  -0.51%  "
$ws.Range("E30").Value = "  -0.63%  "
$ws.Range("E32").Value = "  -1.18%  "
$ws.Range("D33").Value = "1.442.77"
$ws.Range("E33").Value = "  +0.56%  "
$ws.Range("E34").Value = "  -3.36%  "
$ws.Range("E35").Value = "  -4.02%  "
$ws.Range("E36").Value = "  -0.22%  "
$ws.Range("D37").Value = "'0.938"
$ws.Range("E37").Value = "  +4.39%  "
$ws.Range("D38").Value = "'0.561"
$ws.Range("E38").Value = "  -2.19%  "
$ws.Range("E39").Value = "  -0.29%  "
$ws.Range("E40").Value = "  -2.36%  "
$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  +0.23%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "'68.53"
$ws.Range("E42").Value = "  +4.75%  "
$ws.Range("E43").Value = "  -3.20%  "
$ws.Range("E44").Value = "  -0.36%  "
$ws.Range("D45").Value = "'5.42"
$ws.Range("E45").Value = "  -3.07%  "
$ws.Range("E46").Value = "  -2.62%  "
$ws.Range("D47").Value = "1.759.97"
$ws.Range("E47").Value = "  -1.73%  "
$ws.Range("E48").Value = "  +0.07%  "
$ws.Range("D49").Value = "'86.26"
$ws.Range("E49").Value = "  -0.35%  "
$ws.Range("E50").Value = "  +4.97%  "
$ws.Range("D51").Value = "'0.0998"
$ws.Range("E51").Value = "  +0.85%  "
